# Updated cryptos list on Sat May 11 17:12:25 UTC 2024 with GitHub Actions
#
# Applies the scraped-data refresh: new Price (column D) / Volume(1h) (column
# E) figures for (almost) every coin row, plus a rank swap between the
# "Stacks" and "Kaspa" rows (40/41) where name, link, price and volume all
# move together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Many "Price" cells hold strings that *look* numeric (e.g. "6.79", "0.506").
# Excel's COM layer auto-coerces Range.Value assignments of such strings to
# real numbers, same as typing them in interactively. The source data must
# stay text (t="inlineStr" in the XML, no special number formatting), so we
# force text entry with a leading apostrophe and then reset the cell style
# back to "Normal" (the apostrophe/number-format detour otherwise leaves a
# quotePrefix/@ style behind).
function Set-TextCell($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

# --- Row 2: Bitcoin ---
Set-TextCell "D2" "61.110.28"
$ws.Range("E2").Value = "  +0.43%  "

# --- Row 3: Ethereum ---
Set-TextCell "D3" "2.919.90"
$ws.Range("E3").Value = "  +0.40%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  +0.00%  "

# --- Row 5: BNB ---
Set-TextCell "D5" "589.54"
$ws.Range("E5").Value = "  +0.70%  "

# --- Row 6: Solana ---
Set-TextCell "D6" "145.07"
$ws.Range("E6").Value = "  -1.68%  "

# --- Row 7: USDC ---
$ws.Range("E7").Value = "  +0.04%  "

# --- Row 8: XRP ---
Set-TextCell "D8" "0.506"
$ws.Range("E8").Value = "  +0.81%  "

# --- Row 9: LidoStakedEther ---
Set-TextCell "D9" "2.919.70"
$ws.Range("E9").Value = "  +0.32%  "

# --- Row 10: Toncoin ---
Set-TextCell "D10" "6.79"
$ws.Range("E10").Value = "  +1.11%  "

# --- Row 11: Dogecoin ---
$ws.Range("E11").Value = "  +0.20%  "

# --- Row 12: Cardano ---
Set-TextCell "D12" "0.442"
$ws.Range("E12").Value = "  -0.89%  "

# --- Row 13: ShibaInu ---
$ws.Range("E13").Value = "  +1.09%  "

# --- Row 14: Avalanche ---
Set-TextCell "D14" "33.72"
$ws.Range("E14").Value = "  -1.98%  "

# --- Row 15: TRON ---
$ws.Range("E15").Value = "  -0.24%  "

# --- Row 16: WrappedliquidstakedEther2.0 ---
Set-TextCell "D16" "3.405.89"
$ws.Range("E16").Value = "  +0.30%  "

# --- Row 17: WrappedBTC ---
Set-TextCell "D17" "61.062.95"
$ws.Range("E17").Value = "  +0.24%  "

# --- Row 18: Polkadot ---
Set-TextCell "D18" "6.70"
$ws.Range("E18").Value = "  -1.85%  "

# --- Row 19: WrappedEther ---
Set-TextCell "D19" "2.921.97"
$ws.Range("E19").Value = "  +0.22%  "

# --- Row 20: BitcoinCash ---
Set-TextCell "D20" "430.91"
$ws.Range("E20").Value = "  +1.40%  "

# --- Row 21: Chainlink ---
Set-TextCell "D21" "13.45"
$ws.Range("E21").Value = "  -1.46%  "

# --- Row 22: Polygon ---
Set-TextCell "D22" "0.682"
$ws.Range("E22").Value = "  +1.68%  "

# --- Row 23: Uniswap ---
Set-TextCell "D23" "7.08"
$ws.Range("E23").Value = "  -1.08%  "

# --- Row 24: Litecoin ---
Set-TextCell "D24" "80.95"
$ws.Range("E24").Value = "  -0.08%  "

# --- Row 25: RenderToken ---
Set-TextCell "D25" "10.96"
$ws.Range("E25").Value = "  -0.75%  "

# --- Row 26: Fetch.AI ---
Set-TextCell "D26" "2.23"
$ws.Range("E26").Value = "  +2.19%  "

# --- Row 27: InternetComputer(DFINITY) ---
Set-TextCell "D27" "12.03"
$ws.Range("E27").Value = "  +1.80%  "

# --- Row 28: Dai ---
$ws.Range("E28").Value = "  +0.07%  "

# --- Row 29: ImmutableX ---
Set-TextCell "D29" "2.31"
$ws.Range("E29").Value = "  +5.66%  "

# --- Row 30: FirstDigitalUSD ---
$ws.Range("E30").Value = "  -0.06%  "

# --- Row 31: PancakeSwap ---
$ws.Range("E31").Value = "  -0.01%  "

# --- Row 32: NEARProtocol ---
Set-TextCell "D32" "7.13"
$ws.Range("E32").Value = "  -2.33%  "

# --- Row 33: EthereumClassic ---
Set-TextCell "D33" "26.51"
$ws.Range("E33").Value = "  -0.82%  "

# --- Row 34: Hedera ---
$ws.Range("E34").Value = "  +1.64%  "

# --- Row 35: PEPE ---
Set-TextCell "D35" "0.0₃0863"
$ws.Range("E35").Value = "  +2.93%  "

# --- Row 36: Mantle ---
$ws.Range("E36").Value = "  +0.61%  "

# --- Row 37: Filecoin ---
Set-TextCell "D37" "5.63"
$ws.Range("E37").Value = "  -0.76%  "

# --- Row 38: dogwifhat ---
Set-TextCell "D38" "3.07"
$ws.Range("E38").Value = "  +2.70%  "

# --- Row 39: OKB ---
Set-TextCell "D39" "49.92"
$ws.Range("E39").Value = "  +0.43%  "

# --- Rows 40/41: Stacks and Kaspa swap ranking positions ---
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell "D40" "0.125"
$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell "D41" "2.00"
$ws.Range("E41").Value = "  -1.74%  "

# --- Row 42: Cosmos ---
$ws.Range("E42").Value = "  -1.75%  "

# --- Row 43: TheGraph ---
$ws.Range("E43").Value = "  -0.26%  "

# --- Row 44: Arweave ---
Set-TextCell "D44" "39.31"
$ws.Range("E44").Value = "  -5.42%  "

# --- Row 45: Bittensor ---
Set-TextCell "D45" "377.09"
$ws.Range("E45").Value = "  -0.28%  "

# --- Row 46: VeChain ---
$ws.Range("E46").Value = "  +0.71%  "

# --- Row 47: Maker ---
Set-TextCell "D47" "2.709.87"
$ws.Range("E47").Value = "  +2.10%  "

# --- Row 48: Monero ---
Set-TextCell "D48" "131.68"
$ws.Range("E48").Value = "  -1.12%  "

# --- Row 49: USDe --- (unchanged)

# --- Row 50: InjectiveProtocol ---
Set-TextCell "D50" "24.16"
$ws.Range("E50").Value = "  -5.12%  "

# --- Row 51: Stellar ---
$ws.Range("E51").Value = "  +0.24%  "
